$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with rows 461-464 holding four price records.
# This edit inserts four *new* price records right before the current
# row 461, pushing the existing rows 461-464 down to become rows 465-468
# (their contents stay exactly the same, only their row numbers shift).

# Insert 4 blank rows starting at row 461 (shifts old 461:464 -> 465:468)
$ws.Range("A461:T464").EntireRow.Insert()

# ---- New row 461 ----
$ws.Cells.Item(461, 1).Value = 10
$ws.Cells.Item(461, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(461, 3).Value = "La Araucanía"
$ws.Cells.Item(461, 4).Value = 44595
$ws.Cells.Item(461, 5).Value = 9
$ws.Cells.Item(461, 6).Value = "Fruta"
$ws.Cells.Item(461, 7).Value = 100108
$ws.Cells.Item(461, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(461, 9).Value = 100108006
$ws.Cells.Item(461, 10).Value = "Plátano"
$ws.Cells.Item(461, 11).Value = "Barraganete"
$ws.Cells.Item(461, 12).Value = "Primera"
$ws.Cells.Item(461, 13).Value = 55
$ws.Cells.Item(461, 14).Value = 25000
$ws.Cells.Item(461, 15).Value = 25000
$ws.Cells.Item(461, 16).Value = 25000
$ws.Cells.Item(461, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(461, 18).Value = "Ecuador"
$ws.Cells.Item(461, 19).Value = 1250
$ws.Cells.Item(461, 20).Value = 20

# ---- New row 462 ----
$ws.Cells.Item(462, 1).Value = 10
$ws.Cells.Item(462, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(462, 3).Value = "La Araucanía"
$ws.Cells.Item(462, 4).Value = 44595
$ws.Cells.Item(462, 5).Value = 9
$ws.Cells.Item(462, 6).Value = "Fruta"
$ws.Cells.Item(462, 7).Value = 100108
$ws.Cells.Item(462, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(462, 9).Value = 100108006
$ws.Cells.Item(462, 10).Value = "Plátano"
$ws.Cells.Item(462, 11).Value = "Sin especificar"
$ws.Cells.Item(462, 12).Value = "Maduro"
$ws.Cells.Item(462, 13).Value = 125
$ws.Cells.Item(462, 14).Value = 13000
$ws.Cells.Item(462, 15).Value = 13000
$ws.Cells.Item(462, 16).Value = 13000
$ws.Cells.Item(462, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(462, 18).Value = "Ecuador"
$ws.Cells.Item(462, 19).Value = 650
$ws.Cells.Item(462, 20).Value = 20

# ---- New row 463 ----
$ws.Cells.Item(463, 1).Value = 10
$ws.Cells.Item(463, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(463, 3).Value = "La Araucanía"
$ws.Cells.Item(463, 4).Value = 44595
$ws.Cells.Item(463, 5).Value = 9
$ws.Cells.Item(463, 6).Value = "Fruta"
$ws.Cells.Item(463, 7).Value = 100108
$ws.Cells.Item(463, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(463, 9).Value = 100108006
$ws.Cells.Item(463, 10).Value = "Plátano"
$ws.Cells.Item(463, 11).Value = "Sin especificar"
$ws.Cells.Item(463, 12).Value = "Pintón"
$ws.Cells.Item(463, 13).Value = 1250
$ws.Cells.Item(463, 14).Value = 16000
$ws.Cells.Item(463, 15).Value = 17000
$ws.Cells.Item(463, 16).Value = 16520
$ws.Cells.Item(463, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(463, 18).Value = "Ecuador"
$ws.Cells.Item(463, 19).Value = 826
$ws.Cells.Item(463, 20).Value = 20

# ---- New row 464 ----
$ws.Cells.Item(464, 1).Value = 10
$ws.Cells.Item(464, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(464, 3).Value = "La Araucanía"
$ws.Cells.Item(464, 4).Value = 44595
$ws.Cells.Item(464, 5).Value = 9
$ws.Cells.Item(464, 6).Value = "Fruta"
$ws.Cells.Item(464, 7).Value = 100108
$ws.Cells.Item(464, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(464, 9).Value = 100108006
$ws.Cells.Item(464, 10).Value = "Plátano"
$ws.Cells.Item(464, 11).Value = "Sin especificar"
$ws.Cells.Item(464, 12).Value = "Verde"
$ws.Cells.Item(464, 13).Value = 125
$ws.Cells.Item(464, 14).Value = 18000
$ws.Cells.Item(464, 15).Value = 18000
$ws.Cells.Item(464, 16).Value = 18000
$ws.Cells.Item(464, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(464, 18).Value = "Ecuador"
$ws.Cells.Item(464, 19).Value = 900
$ws.Cells.Item(464, 20).Value = 20
